$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '96.407.62'
$ws.Range('E2').Value = '  -1.20%  '
$ws.Range('D3').Value = '3.328.48'
$ws.Range('E3').Value = '  -2.47%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = "'248.85"
$ws.Range('E5').Value = '  -2.27%  '
$ws.Range('D6').Value = "'652.08"
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('E7').Value = '  -6.35%  '
$ws.Range('D8').Value = "'0.419"
$ws.Range('E8').Value = '  -1.15%  '
$ws.Range('D9').Value = "'0.999"
$ws.Range('D10').Value = "'0.986"
$ws.Range('E10').Value = '  -7.47%  '
$ws.Range('D11').Value = '3.325.17'
$ws.Range('E11').Value = '  -2.47%  '
$ws.Range('E12').Value = '  -3.35%  '
$ws.Range('D13').Value = "'40.26"
$ws.Range('E13').Value = '  -3.79%  '
$ws.Range('D14').Value = '96.090.01'
$ws.Range('E14').Value = '  -1.21%  '
$ws.Range('D15').Value = "'6.06"
$ws.Range('E15').Value = '  -3.42%  '
$ws.Range('D16').Value = "'0.0000250"
$ws.Range('E16').Value = '  -3.78%  '
$ws.Range('D17').Value = '3.938.52'
$ws.Range('E17').Value = '  -2.67%  '
$ws.Range('D18').Value = "'8.48"
$ws.Range('E18').Value = '  -1.69%  '
$ws.Range('D19').Value = '3.334.87'
$ws.Range('E19').Value = '  -2.30%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').Value = "'16.99"
$ws.Range('E20').Value = '  -3.13%  '
$ws.Range('B21').Value = 'Stellar'
$ws.Range('C21').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D21').Value = "'0.523"
$ws.Range('E21').Value = '  +4.02%  '
$ws.Range('D22').Value = "'502.97"
$ws.Range('E22').Value = '  -0.71%  '
$ws.Range('D23').Value = "'3.37"
$ws.Range('E23').Value = '  -1.99%  '
$ws.Range('D24').Value = "'10.45"
$ws.Range('D25').Value = "'0.0000197"
$ws.Range('E25').Value = '  -3.76%  '
$ws.Range('D26').Value = "'6.54"
$ws.Range('E26').Value = '  +6.40%  '
$ws.Range('D27').Value = "'95.56"
$ws.Range('E27').Value = '  -3.27%  '
$ws.Range('D28').Value = "'12.01"
$ws.Range('E28').Value = '  -5.87%  '
$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').Value = "'0.143"
$ws.Range('E29').Value = '  -8.41%  '
$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').Value = "'1.00"
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = "'10.96"
$ws.Range('E31').Value = '  -4.03%  '
$ws.Range('B32').Value = 'Cronos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D32').Value = "'0.188"
$ws.Range('E32').Value = '  -5.60%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = "'2.46"
$ws.Range('E33').Value = '  +8.41%  '
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').Value = "'0.998"
$ws.Range('E34').Value = '  -0.32%  '
$ws.Range('B35').Value = 'PolygonEcosystemToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D35').Value = "'0.544"
$ws.Range('E35').Value = '  -5.51%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').Value = "'27.87"
$ws.Range('E36').Value = '  -6.86%  '
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D37').Value = "'1.46"
$ws.Range('E37').Value = '  +3.57%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D38').Value = "'7.59"
$ws.Range('E38').Value = '  -1.93%  '
$ws.Range('B39').Value = 'USDe'
$ws.Range('C39').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D39').Value = "'1.00"
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Value = "'0.150"
$ws.Range('E40').Value = '  -2.24%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').Value = "'505.67"
$ws.Range('E41').Value = '  -1.34%  '
$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').Value = "'24.33"
$ws.Range('E42').Value = '  -1.53%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = "'0.0427"
$ws.Range('E43').Value = '  +1.22%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').Value = "'0.827"
$ws.Range('E44').Value = '  -3.66%  '
$ws.Range('B45').Value = 'MantraDAO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D45').Value = "'3.63"
$ws.Range('E45').Value = '  -0.55%  '
$ws.Range('B46').Value = 'ImmutableX'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D46').Value = "'1.66"
$ws.Range('E46').Value = '  +5.66%  '
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').Value = "'5.45"
$ws.Range('E47').Value = '  -1.09%  '
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').Value = "'8.31"
$ws.Range('E48').Value = '  +1.28%  '
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').Value = "'53.44"
$ws.Range('E49').Value = '  +3.53%  '
$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D50').Value = "'3.10"
$ws.Range('E50').Value = '  -5.06%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').Value = "'162.25"
$ws.Range('E51').Value = '  +0.79%  '
